# Auto-generated edit script: apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.935.98'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '3.848.84'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '697.98'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.07'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").Value = '3.847.13'
$ws.Range("E7").Value = '  +1.33%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.523'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.27'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.17'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").Value = '4.498.87'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '3.852.69'
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("D17").Value = '70.979.36'
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("E20").Value = '  -3.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.73'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.94'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.718'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.69'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000147'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.20'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.65%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.57'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("E28").Value = '  -2.87%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.12'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.51'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.35%  '
$ws.Range("E32").Value = '  -2.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.47'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.61%  '
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("D35").Value = '3.807.43'
$ws.Range("E35").Value = '  +1.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.16'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("E39").Value = '  +6.27%  '
$ws.Range("E40").Value = '  +7.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.00'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.34'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.43%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '163.83'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.83%  '
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000310'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.77'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.47%  '
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.62'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.17'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '408.49'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.78%  '
